$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.769.99'
$ws.Range("E2").Value = '  -1.85%  '
$ws.Range("D3").Value = '1.548.56'
$ws.Range("E3").Value = '  -1.74%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '204.80'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.51%  '
$ws.Range("E6").Value = '  -1.75%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("B8").Value = 'Solana'
$ws.Range("C8").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '21.36'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.11%  '
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.245'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.17%  '
$ws.Range("E10").Value = '  -1.78%  '
$ws.Range("E11").Value = '  -0.90%  '
$ws.Range("E12").Value = '  -1.76%  '
$ws.Range("D13").Value = '1.547.91'
$ws.Range("E13").Value = '  -1.80%  '
$ws.Range("E14").Value = '  -2.74%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.511'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.79%  '
$ws.Range("D16").Value = '26.755.44'
$ws.Range("E16").Value = '  -1.99%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '60.87'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.56%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '213.49'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.75%  '
$ws.Range("E20").Value = '  -1.59%  '
$ws.Range("E22").Value = '  -1.68%  '
$ws.Range("E23").Value = '  -4.26%  '
$ws.Range("E24").Value = '  -0.92%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.78'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.56%  '
$ws.Range("E26").Value = '  -2.54%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.86'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.73%  '
$ws.Range("E28").Value = '  -0.01%  '
$ws.Range("E29").Value = '  -2.00%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0461'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.96%  '
$ws.Range("E31").Value = '  -3.84%  '
$ws.Range("E32").Value = '  -0.32%  '
$ws.Range("D33").Value = '1.349.82'
$ws.Range("E33").Value = '  -4.29%  '
$ws.Range("E34").Value = '  -0.81%  '
$ws.Range("E35").Value = '  -3.71%  '
$ws.Range("E36").Value = '  -0.83%  '
$ws.Range("E37").Value = '  -2.31%  '
$ws.Range("E38").Value = '  -2.15%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.522'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.60%  '
$ws.Range("E40").Value = '  -2.35%  '
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.990'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.23%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.55'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.94%  '
$ws.Range("E44").Value = '  +0.21%  '
$ws.Range("E45").Value = '  -3.26%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '62.87'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.69%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.26'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.63%  '
$ws.Range("D48").Value = '1.682.33'
$ws.Range("E48").Value = '  -1.87%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '85.82'
$ws.Range("D49").Style = "Normal"
$ws.Range("E50").Value = '  +2.38%  '
$ws.Range("D51").Value = '0.0₇0972'
$ws.Range("E51").Value = '  -1.35%  '
